# Extend Report and Custom DataProvider
# Adds a new worksheet "OpenAccountTest" with customer/currency data,
# and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "OpenAccountTest"

# Populate header row.
$newSheet.Range("A1").Value = "customer"
$newSheet.Range("B1").Value = "currency"

# Populate data row.
$newSheet.Range("A2").Value = "Alice Johnson"
$newSheet.Range("B2").Value = "Dollar"

# Autosize column A to fit its contents (matches the "bestFit" column
# width Excel stores after the data is entered).
$newSheet.Columns.Item(1).AutoFit() | Out-Null

# Select B2 like in the target worksheet.
$newSheet.Range("B2").Select()

# Make the new sheet the active one (tab 2 / index 1 zero-based).
$newSheet.Activate()
